# Made fixes for Prod/Demo Verification Script and Display Convenience fees
# qa done for all versions.
#
# Updates the "DateProd" (column B) verification timestamps on each of the
# VT-* sheets to reflect the latest QA pass (Wed May 21 2025 run).

$wb = $excel.ActiveWorkbook

$updates = @{
    "VT-SaleVoid-DualCF-Generic" = @{
        "B2" = "Wed May 21 15:35:56 IST 2025"
        "B3" = "Wed May 21 15:36:45 IST 2025"
        "B4" = "Wed May 21 15:37:23 IST 2025"
        "B5" = "Wed May 21 15:38:13 IST 2025"
    }
    "VT-SaleVoid-NoCF-Generic" = @{
        "B2" = "Wed May 21 15:38:51 IST 2025"
        "B3" = "Wed May 21 15:39:32 IST 2025"
        "B4" = "Wed May 21 15:40:22 IST 2025"
        "B5" = "Wed May 21 15:41:25 IST 2025"
    }
    "VT-SaleVoid-SingleCF-Generic" = @{
        "B2" = "Wed May 21 15:42:25 IST 2025"
        "B3" = "Wed May 21 15:43:20 IST 2025"
        "B4" = "Wed May 21 15:44:03 IST 2025"
        "B5" = "Wed May 21 15:44:58 IST 2025"
    }
    "VT-SaleCredit-DualCF-Generic" = @{
        "B2" = "Wed May 21 15:29:27 IST 2025"
        "B3" = "Wed May 21 15:30:11 IST 2025"
        "B4" = "Wed May 21 15:30:54 IST 2025"
        "B5" = "Wed May 21 15:31:30 IST 2025"
    }
    "VT-SaleCredit-SingleCF-Generic" = @{
        "B2" = "Wed May 21 15:32:15 IST 2025"
        "B3" = "Wed May 21 15:33:18 IST 2025"
        "B4" = "Wed May 21 15:34:20 IST 2025"
        "B5" = "Wed May 21 15:35:18 IST 2025"
    }
    "VT-AuthCapCredit-Generic" = @{
        "B2" = "Wed May 21 14:43:17 IST 2025"
        "B3" = "Wed May 21 14:44:39 IST 2025"
        "B4" = "Wed May 21 14:45:55 IST 2025"
        "B5" = "Wed May 21 14:47:25 IST 2025"
        "B6" = "Wed May 21 14:48:43 IST 2025"
        "B7" = "Wed May 21 14:50:06 IST 2025"
    }
    "VT-AuthCapVoid-Generic" = @{
        "B2" = "Wed May 21 14:51:22 IST 2025"
        "B3" = "Wed May 21 14:52:08 IST 2025"
        "B4" = "Wed May 21 14:52:53 IST 2025"
        "B5" = "Wed May 21 14:53:59 IST 2025"
        "B6" = "Wed May 21 14:55:27 IST 2025"
        "B7" = "Wed May 21 14:56:24 IST 2025"
    }
    "VT-ManualAuthCapture-Generic" = @{
        "B2" = "Wed May 21 15:24:16 IST 2025"
        "B3" = "Wed May 21 15:25:03 IST 2025"
        "B4" = "Wed May 21 15:25:54 IST 2025"
        "B5" = "Wed May 21 15:26:50 IST 2025"
        "B6" = "Wed May 21 15:27:53 IST 2025"
        "B7" = "Wed May 21 15:28:44 IST 2025"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
